$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("APITestData")

$body = "{`n  `"userId`": `"string`",`n  `"collectionOfIsbns`": [`n    {`n      `"isbn`": `"9781449325862`"`n    }`n  ]`n}"

$ws.Range("A7").Value = "AddBooks"
$ws.Range("B7").Value = "/BookStore/v1/Books"
$ws.Range("C7").Value = "Authorization"
$ws.Range("D7").Value = "Bearer"
$ws.Range("G7").Value = $body

$ws.Range("G7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 174

$ws.Range("G3").Select()
